# Apply requested updates:
# 1. Insert a new "gene" column before the existing "chr_pos" column (H),
#    pushing chr_pos to column I.
# 2. Populate the new "gene" column (H2:H17) with values cycling A1/B1/C1/D1
#    in step with the chr_pos values already present in each row.
# 3. Correct a handful of refDepth/altDepth metric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at H, which shifts the old H ("chr_pos") column to I,
# and leaves a blank column at H for the new "gene" field.
$ws.Columns("H:H").Insert()

# New header for column H.
$ws.Range("H1").Value = "gene"

# Gene values per row, aligned with the chr_pos values already in column I.
$gene = @{
    2  = "A1"
    3  = "B1"
    4  = "C1"
    5  = "D1"
    6  = "A1"
    7  = "B1"
    8  = "C1"
    9  = "D1"
    10 = "A1"
    11 = "B1"
    12 = "C1"
    13 = "D1"
    14 = "A1"
    15 = "B1"
    16 = "C1"
    17 = "D1"
}

foreach ($row in $gene.Keys) {
    $ws.Cells.Item($row, 8).Value = $gene[$row]
}

# Corrected metric values.
$ws.Range("E4").Value = 0
$ws.Range("D8").Value = 48
$ws.Range("E11").Value = 0
$ws.Range("D15").Value = 49
$ws.Range("D16").Value = 48
